$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New Cypher query text for each shared query column/tab ----
$statQuery = 'MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)' + "`r`n" + 'OPTIONAL MATCH (samp:sample)-->(c)' + "`r`n" + 'OPTIONAL MATCH (diag:diagnosis)-->(c)' + "`r`n" + 'OPTIONAL MATCH (f:file)-[*]->(c)' + "`r`n" + 'OPTIONAL MATCH (sf:file)-->(s)' + "`r`n" + 'WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p' + "`r`n" + 'MATCH (r:registration)-->(c)' + "`r`n" + 'WHERE r.registration_origin = "PURDUE"' + "`r`n" + 'RETURN  ' + "`r`n" + '    count(distinct p) AS Programs,' + "`r`n" + '    count(distinct s) AS Studies,' + "`r`n" + '    count(distinct c) AS Cases,' + "`r`n" + '    count(distinct samp) AS Samples,' + "`r`n" + '    count(distinct f) AS `Case Files`,' + "`r`n" + '    count(distinct sf) AS `Study Files`'
$casesQuery = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis), (r:registration)-->(c)' + "`r`n" + 'WHERE r.registration_origin = "PURDUE"' + "`r`n" + 'MATCH (c)<--(diag:diagnosis)' + "`r`n" + 'OPTIONAL MATCH (samp:sample)-->(c)' + "`r`n" + 'OPTIONAL MATCH (co:cohort)<-[*]-(c)' + "`r`n" + 'WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight' + "`r`n" + 'RETURN  ' + "`r`n" + '       coalesce(c.case_id, '''') AS `Case ID`,' + "`r`n" + '       coalesce(s.clinical_study_designation, '''') AS `Study Code`,' + "`r`n" + '       coalesce(s.clinical_study_type, '''') AS  `Study Type`,' + "`r`n" + '       coalesce(demo.breed, '''') AS Breed ,' + "`r`n" + '       coalesce(diag.disease_term, '''') AS Diagnosis ,' + "`r`n" + '       coalesce(diag.stage_of_disease, '''') AS `Stage of Disease`,' + "`r`n" + '       CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END AS Age,' + "`r`n" + '       coalesce(demo.sex, '''') AS Sex,' + "`r`n" + '       coalesce(demo.neutered_indicator, '''') AS `Neutered Status`,' + "`r`n" + '       coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '''') AS `Weight (kg)`,' + "`r`n" + '       coalesce(diag.best_response, '''') AS `Response to Treatment`,' + "`r`n" + '       coalesce(co.cohort_description, '''') AS `Cohort`' + "`r`n" + 'Order by c.case_id LIMIT 100'
$samplesQuery = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis), (r:registration)-->(c) ' + "`r`n" + 'WHERE r.registration_origin = "PURDUE"' + "`r`n" + 'WITH DISTINCT samp AS samp, c, demo, diag' + "`r`n" + 'RETURN  coalesce(samp.sample_id, '''') AS `Sample ID`, ' + "`r`n" + '        coalesce(c.case_id, '''') AS `Case ID`, ' + "`r`n" + '        coalesce(demo.breed,'''') AS Breed,' + "`r`n" + '        coalesce(diag.disease_term,'''') AS Diagnosis, ' + "`r`n" + '        coalesce(samp.sample_site, '''') AS `Sample Site`,' + "`r`n" + '        coalesce(samp.summarized_sample_type, '''') AS `Sample Type`,' + "`r`n" + '        coalesce(samp.specific_sample_pathology, '''') AS `Pathology/Morphology`,' + "`r`n" + '        coalesce(samp.sample_chronology, '''') AS `Sample Chronology`,' + "`r`n" + '        coalesce(samp.percentage_tumor, '''') AS `Percentage Tumor`,' + "`r`n" + '        coalesce(samp.necropsy_sample, '''') AS `Necropsy Sample`,' + "`r`n" + '        coalesce(samp.sample_preservation, '''') AS `Sample Preservation`' + "`r`n" + 'Order by samp.sample_id LIMIT 100'
$filesQuery = 'MATCH (f:file)-->(parent)' + "`r`n" + 'WITH DISTINCT f, parent' + "`r`n" + 'MATCH (diag:diagnosis)-->(c)' + "`r`n" + 'OPTIONAL MATCH (f)-[*]->(samp:sample)' + "`r`n" + 'MATCH (f)-[*]->(c:case)<--(demo:demographic)' + "`r`n" + 'MATCH (r:registration)-->(c)' + "`r`n" + 'WHERE r.registration_origin = "PURDUE"' + "`r`n" + 'OPTIONAL MATCH (s:study)<--(c)<--(diag:diagnosis)<-[*]-(samp)' + "`r`n" + 'WITH' + "`r`n" + '        f, parent, c, demo, diag, s, samp,' + "`r`n" + '        [''Bytes'', ''KB'', ''MB'', ''GB'', ''TB''] AS units,' + "`r`n" + '        toInteger(floor(log(f.file_size)/log(1024))) as i,' + "`r`n" + '        2 as precision' + "`r`n" + 'WITH' + "`r`n" + '        f, parent, c, demo, diag, s, samp,' + "`r`n" + '        f.file_size /(1024^i) AS value, ' + "`r`n" + '        10^precision AS factor,' + "`r`n" + '        units[i] as unit' + "`r`n" + 'WITH    ' + "`r`n" + '        f, parent, c, demo, diag, s, samp, unit,' + "`r`n" + '        round(factor * value)/factor AS size' + "`r`n" + 'RETURN ' + "`r`n" + '        coalesce(f.file_name, '''') AS `File Name`,' + "`r`n" + '        coalesce(f.file_format, '''') AS `Format`,' + "`r`n" + '        coalesce(f.file_type, '''') AS `File Type`,' + "`r`n" + '       CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+'' '' +unit ELSE size+'' '' +unit END AS Size,' + "`r`n" + '        coalesce(labels(parent)[0], '''') AS `Association`,' + "`r`n" + '        coalesce(f.file_description, '''') AS `Description`,' + "`r`n" + '        coalesce(samp.sample_id, '''') AS `Sample ID`,' + "`r`n" + '        coalesce(c.case_id, '''') AS `Case ID`,' + "`r`n" + '        coalesce(demo.breed,'''') AS Breed ,' + "`r`n" + '        coalesce(diag.disease_term,'''') AS Diagnosis' + "`r`n" + 'Order By f.file_name LIMIT 100'
$studyFilesQuery = 'MATCH (f:file)-->(s:study)' + "`r`n" + 'MATCH (s)<--(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)' + "`r`n" + 'MATCH (r:registration)-->(c)' + "`r`n" + 'WHERE r.registration_origin = "PURDUE"' + "`r`n" + 'WITH DISTINCT f,  s, c, demo, diag' + "`r`n" + 'WITH' + "`r`n" + '        f, c, demo, diag, s,' + "`r`n" + '        [''Bytes'', ''KB'', ''MB'', ''GB'', ''TB''] AS units,' + "`r`n" + '        toInteger(floor(log(f.file_size)/log(1024))) as i,' + "`r`n" + '        2 as precision' + "`r`n" + 'WITH    ' + "`r`n" + '        f, c, demo, diag, s,' + "`r`n" + '        f.file_size /(1024^i) AS value, 10^precision AS factor,' + "`r`n" + '        units[i] as unit' + "`r`n" + '        WITH    ' + "`r`n" + '        f,  c, demo, diag, s, unit,' + "`r`n" + '        round(factor * value)/factor AS size' + "`r`n" + 'RETURN DISTINCT' + "`r`n" + '  coalesce(f.file_name, '''') AS `File Name`,' + "`r`n" + '  coalesce(f.file_type, '''') AS `File Type`,' + "`r`n" + '  coalesce("study", '''') AS `Association`,' + "`r`n" + '  coalesce(f.file_description, '''') AS `Description`,' + "`r`n" + '  coalesce(f.file_format, '''') AS  Format,' + "`r`n" + '  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+'' '' +unit ELSE size+'' '' +unit END AS Size,' + "`r`n" + '  coalesce(s.clinical_study_designation,'''') AS `Study Code`' + "`r`n" + '  order by ''File Name'' asc' + "`r`n" + '  limit 100'

# ---- StatQuery (column C) is shared by every tab row ----
$ws.Range("C2").Value = $statQuery
$ws.Range("C3").Value = $statQuery
$ws.Range("C4").Value = $statQuery
$ws.Range("C5").Value = $statQuery

# ---- CasesTab (row 2) ----
$ws.Range("B2").Value = $casesQuery
$ws.Rows(2).RowHeight = 333.5

# ---- SamplesTab (row 3) ----
$ws.Range("B3").Value = $samplesQuery
$ws.Rows(3).RowHeight = 246.5

# ---- FilesTab (row 4) ----
$ws.Range("B4").Value = $filesQuery
$ws.Rows(4).RowHeight = 409.5

# ---- StudyFilesTab (row 5) ----
$ws.Range("B5").Value = $studyFilesQuery
$ws.Rows(5).RowHeight = 409.5

# ---- Move the active selection to B5 (matches the author's final cursor position) ----
$ws.Range("B5").Select()

